$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B15 was mistakenly entered as text "4" - fix it to be a real number 4
$ws.Range("B15").Value = 4

# Append a new annotation row (row 16) for Sunsi Wu
$ws.Range("A16").Value = "Sunsi Wu"

# politeness_score "3" must stay textual (matches the source export, which
# wrote it as an inline string) - use a leading apostrophe to force text
# entry, then reset the style so no extra number-format style is left behind.
$ws.Range("B16").Value = "'3"
$ws.Range("B16").Style = "Normal"

$ws.Range("C16").Value = "does not; "
$ws.Range("D16").Value = "DFT"
$ws.Range("E16").Value = "WRI"
$ws.Range("F16").Value = "4cbdf296-0ef7-4a60-9d08-bf70fb941ab3"
$ws.Range("G16").Value = "SJTB5GZCb_annotated.xlsx"
$ws.Range("H16").Value = "The paper does not sufficiently discuss and compare the relevant neuroscience literature and related work."
